$wb = $excel.ActiveWorkbook

# --- 1. Insert a new "2022-Q1" sheet right before the "总计" sheet ---------
$totalSheet  = $wb.Worksheets.Item("总计")
$sourceSheet = $wb.Worksheets.Item("2021-Q4")      # donor sheet for header/style layout

$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# NOTE: a worksheet reference obtained via a *positional* index before the
# Add() call resolves positionally, so after inserting a new sheet such a
# handle would now point at the newly-added sheet instead of "总计". Always
# re-fetch "总计" by name after the insertion.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy header row (B1:H1) and the styled index column (A2:A8) from the
# previous quarter sheet so fonts/borders/number formats match exactly.
$sourceSheet.Range("B1:H1").Copy($newSheet.Range("B1"))
$sourceSheet.Range("A2:A8").Copy($newSheet.Range("A2"))

# --- fund holdings data for 2022-Q1 ---
$fundRows = @(
    @("161219", "国投瑞银新兴产业混合(LOF)",     "8.46", "79.83", "5.63", "0.4763", 3),
    @("002628", "招商安博灵活配置混合A",          "1.55", "65.07", "5.94", "0.0921", 2),
    @("161225", "国投瑞银瑞盈灵活配置混合（LOF）", "2.28", "94.46", "3.62", "0.0825", 9),
    @("519097", "新华中小市值优选混合",            "0.75", "62.70", "3.18", "0.0238", 8),
    @("002629", "招商安博灵活配置混合C",          "0.31", "65.07", "5.94", "0.0184", 2),
    @("005169", "华泰保兴策略精选灵活配置混合A",   "0.38", "84.18", "3.27", "0.0124", 9),
    @("005170", "华泰保兴策略精选灵活配置混合C",   "0.23", "84.18", "3.27", "0.0075", 9)
)

$r = 2
foreach ($fr in $fundRows) {
    # B = 基金代码 (fund code) - force Text so leading zeros survive
    $codeCell = $newSheet.Range("B$r")
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $fr[0]

    # C = 基金名称 (fund name) - plain text, no coercion risk
    $newSheet.Range("C$r").Value = $fr[1]

    # D/E/F/G = 基金规模 / 股票总仓位 / 仓位占比 / 持有市值(亿元) - numeric-looking
    # strings that must stay text, as in the source sheets.
    $newSheet.Range("D$r").NumberFormat = "@"
    $newSheet.Range("D$r").Value = $fr[2]
    $newSheet.Range("E$r").NumberFormat = "@"
    $newSheet.Range("E$r").Value = $fr[3]
    $newSheet.Range("F$r").NumberFormat = "@"
    $newSheet.Range("F$r").Value = $fr[4]
    $newSheet.Range("G$r").NumberFormat = "@"
    $newSheet.Range("G$r").Value = $fr[5]

    # H = 仓位排名 (rank) - real number
    $newSheet.Range("H$r").Value = $fr[6]

    $r++
}

# --- 2. Add the 2022-Q1 summary row at the top of the "总计" sheet --------
$totalSheet.Rows("2:2").Insert()

# The inserted row picks up stray formatting from the row above; clear it so
# the new cells land with the default (unstyled) look used by all the other
# data rows on this sheet.
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 0.71

# A2 needs the bold/bordered "index column" look used by A3:A7 - grab it from
# the donor sheet (plain Style assignment doesn't stick, Copy() does).
$sourceSheet.Range("A2").Copy($totalSheet.Range("A2"))

# Renumber the helper index column (A2..A7) sequentially 0..5
for ($i = 0; $i -le 5; $i++) {
    $rowNum = 2 + $i
    $totalSheet.Range("A$rowNum").Value = $i
}
